$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.347.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.157.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.47%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.11%  "

$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0871"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.643.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.06%  "

$ws.Range("E16").Value = "  +9.27%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.159.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "53.275.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.38%  "

$ws.Range("E20").Value = "  +4.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("E23").Value = "  +1.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.29%  "

$ws.Range("E28").Value = "  +2.48%  "

$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("E30").Value = "  +0.05%  "

$ws.Range("E31").Value = "  +2.83%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.39%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.10%  "

$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0499"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.61%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.54%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +10.53%  "

$ws.Range("E40").Value = "  +11.23%  "

$ws.Range("E41").Value = "  +0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.41%  "

$ws.Range("E45").Value = "  +1.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.20%  "

$ws.Range("E47").Value = "  -1.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.088.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +23.35%  "

$ws.Range("E51").Value = "  +4.19%  "
